$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-PlainCell($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

Set-TextCell 'D2' '24.381.33'
Set-PlainCell 'E2' '  -5.49%  '
Set-TextCell 'D3' '1.634.28'
Set-PlainCell 'E3' '  -7.19%  '
Set-TextCell 'D4' '1.004'
Set-PlainCell 'E4' '  +0.01%  '
Set-PlainCell 'E5' '  +0.24%  '
Set-TextCell 'D6' '304.85'
Set-PlainCell 'E6' '  -3.78%  '
Set-TextCell 'D7' '0.3592'
Set-PlainCell 'E7' '  -6.06%  '
Set-TextCell 'D8' '46.79'
Set-PlainCell 'E8' '  -6.86%  '
Set-TextCell 'D9' '0.3229'
Set-PlainCell 'E9' '  -10.48%  '
Set-TextCell 'D10' '1.105'
Set-PlainCell 'E10' '  -9.90%  '
Set-TextCell 'D11' '0.06872'
Set-PlainCell 'E11' '  -10.80%  '
Set-TextCell 'D12' '1.004'
Set-PlainCell 'E12' '  +0.32%  '
Set-TextCell 'D13' '5.906'
Set-PlainCell 'E13' '  -8.75%  '
Set-TextCell 'D14' '19.10'
Set-PlainCell 'E14' '  -11.92%  '
Set-TextCell 'D15' '1.636.26'
Set-PlainCell 'E15' '  -7.18%  '
Set-TextCell 'D16' '6.513'
Set-PlainCell 'E16' '  -8.16%  '
Set-TextCell 'D17' '0.00001042'
Set-PlainCell 'E17' '  -9.88%  '
Set-TextCell 'D18' '0.06506'
Set-PlainCell 'E18' '  -4.09%  '
Set-PlainCell 'E19' '  +0.26%  '
Set-TextCell 'D20' '76.38'
Set-PlainCell 'E20' '  -12.13%  '
Set-PlainCell 'B21' 'Avalanche'
Set-PlainCell 'C21' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 'D21' '15.71'
Set-PlainCell 'E21' '  -11.25%  '
Set-PlainCell 'B22' 'Uniswap'
Set-PlainCell 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D22' '5.878'
Set-PlainCell 'E22' '  -9.89%  '
Set-TextCell 'D23' '11.96'
Set-PlainCell 'E23' '  -8.14%  '
Set-TextCell 'D24' '24.317.23'
Set-PlainCell 'E24' '  -5.62%  '
Set-TextCell 'D25' '2.395'
Set-TextCell 'D26' '2.321'
Set-PlainCell 'E26' '  -20.18%  '
Set-TextCell 'D27' '143.96'
Set-PlainCell 'E27' '  -7.66%  '
Set-TextCell 'D28' '18.54'
Set-PlainCell 'E28' '  -10.98%  '
Set-TextCell 'D29' '1.816.05'
Set-PlainCell 'E29' '  -7.34%  '
Set-TextCell 'D30' '123.70'
Set-PlainCell 'E30' '  -7.55%  '
Set-TextCell 'D31' '1.123'
Set-PlainCell 'E31' '  -7.06%  '
Set-TextCell 'D32' '4.060'
Set-PlainCell 'E32' '  -3.80%  '
Set-TextCell 'D33' '5.640'
Set-PlainCell 'E33' '  -22.13%  '
Set-TextCell 'D34' '0.08347'
Set-PlainCell 'E34' '  -4.72%  '
Set-TextCell 'D35' '1.659'
Set-PlainCell 'E35' '  -8.39%  '
Set-TextCell 'D36' '12.29'
Set-PlainCell 'E36' '  -14.47%  '
Set-TextCell 'D37' '5.098'
Set-PlainCell 'E37' '  -11.16%  '
Set-TextCell 'D38' '0.05974'
Set-TextCell 'D39' '0.02205'
Set-PlainCell 'E39' '  -11.82%  '
Set-PlainCell 'E40' '  -7.85%  '
Set-TextCell 'D41' '8.138'
Set-PlainCell 'E41' '  -13.52%  '
Set-TextCell 'D42' '0.2025'
Set-PlainCell 'E42' '  -10.63%  '
Set-PlainCell 'E43' '  +0.28%  '
Set-TextCell 'D44' '0.5828'
Set-PlainCell 'E44' '  -11.70%  '
Set-TextCell 'D45' '3.710'
Set-PlainCell 'E45' '  -4.86%  '
Set-TextCell 'D46' '12.45'
Set-PlainCell 'E46' '  -13.37%  '
Set-TextCell 'D47' '0.5513'
Set-PlainCell 'E47' '  -13.31%  '
Set-TextCell 'D48' '120.90'
Set-PlainCell 'E48' '  -8.51%  '
Set-TextCell 'D49' '1.916'
Set-PlainCell 'E49' '  -12.00%  '
Set-TextCell 'D50' '0.06878'
Set-PlainCell 'E50' '  -8.35%  '
Set-TextCell 'D51' '73.26'
Set-PlainCell 'E51' '  -9.66%  '
